$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 885; existing rows 885:949 shift down to 886:950
$ws.Rows("885").Insert()

# Populate the newly inserted row 885 with the new record's data
$ws.Range("A885").Value = 3
$ws.Range("B885").Value = "Femacal de La Calera"
$ws.Range("C885").Value = "Coquimbo"
$ws.Range("D885").Value = 45265
$ws.Range("E885").Value = 5
$ws.Range("F885").Value = 100112003
$ws.Range("G885").Value = "Ajo"
$ws.Range("H885").Value = "Chino"
$ws.Range("I885").Value = "Primera"
$ws.Range("J885").Value = 103
$ws.Range("K885").Value = 21000
$ws.Range("L885").Value = 22000
$ws.Range("M885").Value = 21340
$ws.Range("N885").Value = "$/caja 10 kilos"
$ws.Range("O885").Value = "China"
$ws.Range("P885").Value = 2134
$ws.Range("Q885").Value = 10
$ws.Range("R885").Value = "Hortaliza"
